# Mixed Fixture: drop the separate "Away" column and turn the "Home"
# column into a single "Location" column.
#
#   - Column C ("Away") is removed entirely (cells shift left).
#   - The header that used to read "Home" (now column B) is renamed
#     to "Location".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "Away" column (C) - remaining data shifts left.
$ws.Columns("C:C").Delete()

# Rename the surviving column B header from "Home" to "Location".
$ws.Range("B1").Value = "Location"

# Match the author's final selection (B2) shown in the saved file.
$ws.Range("B2").Select() | Out-Null
